$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("qwe", "2025-07-14 17:38:22"),
    @("qwe", "2025-07-14 17:42:46"),
    @("qw", "2025-07-24 21:25:23"),
    @("Reha_Sai", "2025-07-28 15:02:11"),
    @("qw", "2025-07-28 15:26:52"),
    @("qw", "2025-07-30 14:53:52")
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
